# ---------------------------------------------------------------------------
# TestWorkbook_FableExceljs.xlsx -> TestWorkbook_Excel style update
#
# Summary of the edit (per the commit "Update files from TestWorkbook_Excel"):
#   * The "DateTime" column (previously holding the placeholder text "any")
#     on every data sheet now holds real date values (2023-10-14 .. 17),
#     formatted with the built-in short-date number format (mm-dd-yy).
#   * The "Tableless" sheet was recreated (delete + re-add) and repopulated,
#     picking up a fresh sheetId and a 5th data row ("Outer Space").
#   * Selections / active sheet moved around a bit.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "WithTable" sheet: DateTime column C2:C5 -> real dates
# ---------------------------------------------------------------------------
$wsWithTable = $wb.Worksheets.Item("WithTable")

$dates1 = $wsWithTable.Range("C2:C5")
$dates1.NumberFormat = "mm-dd-yy"
$wsWithTable.Range("C2").Value2 = 45213
$wsWithTable.Range("C3").Value2 = 45214
$wsWithTable.Range("C4").Value2 = 45215
$wsWithTable.Range("C5").Value2 = 45216

$wsWithTable.Range("E3").Select()

# ---------------------------------------------------------------------------
# 2. Update the "WithTable_Duplicate" sheet: DateTime column D5:D8 -> dates
# ---------------------------------------------------------------------------
$wsDup = $wb.Worksheets.Item("WithTable_Duplicate")

$dates2 = $wsDup.Range("D5:D8")
$dates2.NumberFormat = "mm-dd-yy"
$wsDup.Range("D5").Value2 = 45213
$wsDup.Range("D6").Value2 = 45214
$wsDup.Range("D7").Value2 = 45215
$wsDup.Range("D8").Value2 = 45216

$wsDup.Range("I26").Select()

# ---------------------------------------------------------------------------
# 3. Recreate the "Tableless" sheet (delete the old one, add a fresh sheet
#    in the same slot) and repopulate it with the same shape of data as
#    "WithTable", but without a Table object.
# ---------------------------------------------------------------------------
$oldTableless = $wb.Worksheets.Item("Tableless")
$oldTableless.Delete()

$wsTableless = $wb.Worksheets.Add($wsDup)
$wsTableless.Name = "Tableless"

$wsTableless.Range("A1").Value = "Numbers"
$wsTableless.Range("B1").Value = "Strings"
$wsTableless.Range("C1").Value = "DateTime"
$wsTableless.Range("D1").Value = "ARCtrl Column"
$wsTableless.Range("E1").Value = "ARCtrl Column "

$wsTableless.Range("A2").Value2 = 1
$wsTableless.Range("B2").Value = "Hello"
$wsTableless.Range("C2").NumberFormat = "mm-dd-yy"
$wsTableless.Range("C2").Value2 = 45213
$wsTableless.Range("D2").Value = "(A) This is part 1 of 2"
$wsTableless.Range("E2").Value = "(A) This is part 2 of 2"

$wsTableless.Range("A3").Value2 = 2
$wsTableless.Range("B3").Value = "World"
$wsTableless.Range("C3").NumberFormat = "mm-dd-yy"
$wsTableless.Range("C3").Value2 = 45214
$wsTableless.Range("E3").Value = "Tests if column names with whitespace at end can be unique"

$wsTableless.Range("A4").Value2 = 3
$wsTableless.Range("B4").Value = "Bye"
$wsTableless.Range("C4").NumberFormat = "mm-dd-yy"
$wsTableless.Range("C4").Value2 = 45215

$wsTableless.Range("A5").Value2 = 4
$wsTableless.Range("B5").Value = "Outer Space"
$wsTableless.Range("C5").NumberFormat = "mm-dd-yy"
$wsTableless.Range("C5").Value2 = 45216

$wsTableless.Columns.Item(3).ColumnWidth = 10.140625

$wsTableless.Range("E2").Select()
$wsTableless.Activate()

# Best-effort: rename the default theme the way newer Excel versions label it.
# (No-op on hosts that don't expose a writable theme name.)
try { $wb.Theme.Name = "Office Theme" } catch {}
